$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (rows 2-11) ---
$ws.Range("F2").Value = 1.01
$ws.Range("R3").Value = 1.16
$ws.Range("F4").Value = 1.7
$ws.Range("I4").Value = 5.6
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 5.1
$ws.Range("L4").Value = 1.29
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 4.1
$ws.Range("P4").Value = 2.08
$ws.Range("Q4").Value = 1.74
$ws.Range("R4").Value = 1.43
$ws.Range("S4").Value = 2.6
$ws.Range("T4").Value = 1.73
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.21
$ws.Range("W4").Value = 2.14
$ws.Range("X4").Value = 22
$ws.Range("Y4").Value = 24
$ws.Range("Z4").Value = 50
$ws.Range("AA4").Value = 140
$ws.Range("AB4").Value = 12
$ws.Range("AC4").Value = 11.5
$ws.Range("AD4").Value = 24
$ws.Range("AE4").Value = 75
$ws.Range("AF4").Value = 14
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 75
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 22
$ws.Range("AL4").Value = 40
$ws.Range("AM4").Value = 110
$ws.Range("AN4").Value = 12
$ws.Range("AO4").Value = 75
$ws.Range("F5").Value = 1.94
$ws.Range("G5").Value = 2.72
$ws.Range("H5").Value = 2.3
$ws.Range("I5").Value = 3.5
$ws.Range("K5").Value = 950
$ws.Range("O5").Value = 1.1
$ws.Range("P5").Value = 2.78
$ws.Range("Q5").Value = 1.37
$ws.Range("R5").Value = 1.84
$ws.Range("T5").Value = 1.4
$ws.Range("U5").Value = 2.7
$ws.Range("V5").Value = 1.41
$ws.Range("X5").Value = 980
$ws.Range("Y5").Value = 980
$ws.Range("Z5").Value = 980
$ws.Range("AA5").Value = 55
$ws.Range("AB5").Value = 980
$ws.Range("AC5").Value = 980
$ws.Range("AD5").Value = 980
$ws.Range("AE5").Value = 980
$ws.Range("AF5").Value = 980
$ws.Range("AG5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("AI5").Value = 980
$ws.Range("AJ5").Value = 980
$ws.Range("AK5").Value = 980
$ws.Range("AL5").Value = 980
$ws.Range("AM5").Value = 980
$ws.Range("AN5").Value = 980
$ws.Range("AO5").Value = 980
$ws.Range("L6").Value = 1.44
$ws.Range("AH6").Value = 20
$ws.Range("H7").Value = 2.86
$ws.Range("I7").Value = 3.15
$ws.Range("N7").Value = 3.15
$ws.Range("F8").Value = 3.6
$ws.Range("H8").Value = 2.42
$ws.Range("I8").Value = 2.68
$ws.Range("J8").Value = 2.76
$ws.Range("O8").Value = 1.67
$ws.Range("P8").Value = 1.41
$ws.Range("J9").Value = 1.09
$ws.Range("H10").Value = 3.7
$ws.Range("L10").Value = 1.54
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 2.64
$ws.Range("O10").Value = 1.51
$ws.Range("Q10").Value = 2.48
$ws.Range("R10").Value = 1.2
$ws.Range("S10").Value = 5
$ws.Range("T10").Value = 2.04
$ws.Range("U10").Value = 1.76
$ws.Range("X10").Value = 980
$ws.Range("Y10").Value = 980
$ws.Range("AA10").Value = 120
$ws.Range("AB10").Value = 980
$ws.Range("AC10").Value = 980
$ws.Range("AE10").Value = 70
$ws.Range("AG10").Value = 980
$ws.Range("AI10").Value = 110
$ws.Range("AL10").Value = 65
$ws.Range("AM10").Value = 210
$ws.Range("AN10").Value = 980
$ws.Range("AO10").Value = 110
$ws.Range("N11").Value = 2.98
$ws.Range("S11").Value = 4.2
$ws.Range("X11").Value = 11.5

# --- Add new rows 13-15 ---
# Row 13
$ws.Range("A13").Value = "Ecuadorian Serie A"
$ws.Range("B13").Value = "'2025-11-10"
$ws.Range("C13").Value = "'21:00:00"
$ws.Range("D13").Value = "Aucas"
$ws.Range("E13").Value = "Delfin"
$ws.Range("F13").Value = 1.5
$ws.Range("G13").Value = 1.68
$ws.Range("H13").Value = 5.8
$ws.Range("I13").Value = 9.4
$ws.Range("J13").Value = 3.65
$ws.Range("K13").Value = 5.6
$ws.Range("L13").Value = 1.01
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 3.8
$ws.Range("O13").Value = 1.28
$ws.Range("P13").Value = 1.97
$ws.Range("Q13").Value = 1.81
$ws.Range("R13").Value = 1.37
$ws.Range("S13").Value = 2.82
$ws.Range("T13").Value = 1.79
$ws.Range("U13").Value = 1.74
$ws.Range("V13").Value = 1.13
$ws.Range("W13").Value = 2.46
$ws.Range("X13").Value = 980
$ws.Range("Y13").Value = 980
$ws.Range("Z13").Value = 75
$ws.Range("AA13").Value = 1000
$ws.Range("AB13").Value = 980
$ws.Range("AC13").Value = 980
$ws.Range("AD13").Value = 980
$ws.Range("AE13").Value = 140
$ws.Range("AF13").Value = 980
$ws.Range("AG13").Value = 980
$ws.Range("AH13").Value = 980
$ws.Range("AI13").Value = 130
$ws.Range("AJ13").Value = 980
$ws.Range("AK13").Value = 980
$ws.Range("AL13").Value = 980
$ws.Range("AM13").Value = 1000
$ws.Range("AN13").Value = 10.5
$ws.Range("AO13").Value = 1000

# Row 14
$ws.Range("A14").Value = "Argentinian Primera Division"
$ws.Range("B14").Value = "'2025-11-10"
$ws.Range("C14").Value = "'21:15:00"
$ws.Range("D14").Value = "Independiente Rivadavia"
$ws.Range("E14").Value = "Central Cordoba (SdE)"
$ws.Range("F14").Value = 2.54
$ws.Range("G14").Value = 2.8
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 3.6
$ws.Range("J14").Value = 2.86
$ws.Range("K14").Value = 3.2
$ws.Range("L14").Value = 1.6
$ws.Range("M14").Value = 1.13
$ws.Range("N14").Value = 2.44
$ws.Range("O14").Value = 1.58
$ws.Range("P14").Value = 1.49
$ws.Range("Q14").Value = 2.74
$ws.Range("R14").Value = 1.17
$ws.Range("S14").Value = 5.7
$ws.Range("T14").Value = 2.16
$ws.Range("U14").Value = 1.73
$ws.Range("V14").Value = 1.39
$ws.Range("W14").Value = 1.56
$ws.Range("X14").Value = 980
$ws.Range("Y14").Value = 10.5
$ws.Range("Z14").Value = 980
$ws.Range("AA14").Value = 75
$ws.Range("AB14").Value = 980
$ws.Range("AC14").Value = 980
$ws.Range("AD14").Value = 980
$ws.Range("AE14").Value = 60
$ws.Range("AF14").Value = 980
$ws.Range("AG14").Value = 980
$ws.Range("AH14").Value = 980
$ws.Range("AI14").Value = 110
$ws.Range("AJ14").Value = 980
$ws.Range("AK14").Value = 980
$ws.Range("AL14").Value = 75
$ws.Range("AM14").Value = 250
$ws.Range("AN14").Value = 980
$ws.Range("AO14").Value = 85

# Row 15
$ws.Range("A15").Value = "Argentinian Primera Division"
$ws.Range("B15").Value = "'2025-11-10"
$ws.Range("C15").Value = "'21:15:00"
$ws.Range("D15").Value = "Argentinos Juniors"
$ws.Range("E15").Value = "Belgrano"
$ws.Range("F15").Value = 1.69
$ws.Range("G15").Value = 1.79
$ws.Range("H15").Value = 6
$ws.Range("I15").Value = 7.4
$ws.Range("J15").Value = 3.55
$ws.Range("K15").Value = 3.85
$ws.Range("L15").Value = 1.5
$ws.Range("M15").Value = 1.1
$ws.Range("N15").Value = 2.8
$ws.Range("O15").Value = 1.46
$ws.Range("P15").Value = 1.61
$ws.Range("Q15").Value = 2.32
$ws.Range("R15").Value = 1.22
$ws.Range("S15").Value = 4.6
$ws.Range("T15").Value = 2.2
$ws.Range("U15").Value = 1.68
$ws.Range("V15").Value = 1.15
$ws.Range("W15").Value = 2.26
$ws.Range("X15").Value = 10.5
$ws.Range("Y15").Value = 21
$ws.Range("Z15").Value = 65
$ws.Range("AA15").Value = 300
$ws.Range("AB15").Value = 6.6
$ws.Range("AC15").Value = 8.6
$ws.Range("AD15").Value = 980
$ws.Range("AE15").Value = 170
$ws.Range("AF15").Value = 11
$ws.Range("AG15").Value = 11.5
$ws.Range("AH15").Value = 980
$ws.Range("AI15").Value = 170
$ws.Range("AJ15").Value = 980
$ws.Range("AK15").Value = 980
$ws.Range("AL15").Value = 65
$ws.Range("AM15").Value = 260
$ws.Range("AN15").Value = 980
$ws.Range("AO15").Value = 290

